$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-8 are being reordered (same row contents, shuffled to new rows).
# Stage the original rows 2-8 into a scratch area first (rows 200-206),
# then copy them back into their new destination rows, so that cell
# formatting/type (inline/shared string) is preserved exactly via Copy,
# rather than re-typing values (which would coerce numeric-looking text
# like "2" into a real number).

for ($i = 0; $i -le 6; $i++) {
    $srcRow = 2 + $i
    $scratchRow = 200 + $i
    $ws.Range("A" + $srcRow + ":D" + $srcRow).Copy($ws.Range("A" + $scratchRow + ":D" + $scratchRow))
}

# destination row -> original (source) row whose content it should now hold
$destToSource = @{
    2 = 6
    3 = 4
    4 = 5
    5 = 8
    6 = 2
    7 = 3
    8 = 7
}

foreach ($destRow in 2..8) {
    $srcRow = $destToSource[$destRow]
    $scratchRow = 200 + ($srcRow - 2)
    $ws.Range("A" + $scratchRow + ":D" + $scratchRow).Copy($ws.Range("A" + $destRow + ":D" + $destRow))
}

# Clean up the scratch rows
$ws.Range("A200:D206").Clear()
